$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.274
    3  = 0.188
    4  = 0.376
    5  = 0.355
    6  = 1.155
    7  = 1.089
    8  = 1.221
    9  = 2.21
    10 = 4.605
    11 = 7.154
    12 = 21.545
    13 = 9.744
    14 = 13.981
    15 = 4.564
    16 = 184.271
    17 = 77.62299999999999
    18 = 78.879
    19 = 17.424
    20 = 2887.745
    21 = 1932.174
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}
